$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Benchmark")

# --- Fix E7/E8 number formats to match E5/E6 (numFmtId 2, "0.00") ---
$ws.Range("E5:E6").Copy() | Out-Null
$ws.Range("E7:E8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Add new benchmark rows 9-12 (Octree radius search results) ---
# Copy formatting from row 8 (same date/category layout) down to the new rows
$ws.Range("A8:H8").Copy() | Out-Null
$ws.Range("A9:H12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 9: Radius Search Demeaned (KdTree, OMP) @ scale 2
$ws.Range("A9").Value2 = 44460
$ws.Range("B9").Value2 = "Geomech Desktop"
$ws.Range("C9").Value2 = 79057628
$ws.Range("D9").Value2 = "Radius Search Demeaned (KdTree, OMP)"
$ws.Range("E9").Value2 = 2
$ws.Range("F9").Value2 = 7540197.516
$ws.Range("G9").Formula = "=F9/1000"
$ws.Range("H9").Formula = "=G9/60"

# Row 10: Radius Search Demeaned (Octree, OMP) @ scale 0.25
$ws.Range("A10").Value2 = 44460
$ws.Range("B10").Value2 = "Geomech Desktop"
$ws.Range("C10").Value2 = 79057628
$ws.Range("D10").Value2 = "Radius Search Demeaned (Octree, OMP)"
$ws.Range("E10").Value2 = 0.25
$ws.Range("F10").Value2 = 74567.8539
$ws.Range("G10").Formula = "=F10/1000"
$ws.Range("H10").Formula = "=G10/60"

# Row 11: Radius Search Demeaned (Octree, OMP) @ scale 2
$ws.Range("A11").Value2 = 44460
$ws.Range("B11").Value2 = "Geomech Desktop"
$ws.Range("C11").Value2 = 79057628
$ws.Range("D11").Value2 = "Radius Search Demeaned (Octree, OMP)"
$ws.Range("E11").Value2 = 2
$ws.Range("F11").Value2 = 2438933.4924
$ws.Range("G11").Formula = "=F11/1000"
$ws.Range("H11").Formula = "=G11/60"

# Row 12: Radius Search Demeaned (Octree, OMP) @ scale 1
$ws.Range("A12").Value2 = 44460
$ws.Range("B12").Value2 = "Geomech Desktop"
$ws.Range("C12").Value2 = 79057628
$ws.Range("D12").Value2 = "Radius Search Demeaned (Octree, OMP)"
$ws.Range("E12").Value2 = 1
$ws.Range("F12").Value2 = 617709.6166
$ws.Range("G12").Formula = "=F12/1000"
$ws.Range("H12").Formula = "=G12/60"

# --- Update the active selection shown in the saved workbook ---
$ws.Range("H10").Select() | Out-Null

Write-Host "done"
